$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of bare channel name -> sampling rate (Hz), read via pyedflib
$rates = @{
    "E1"       = 128
    "E2"       = 128
    "F3"       = 128
    "F4"       = 128
    "C3"       = 128
    "C4"       = 128
    "O1"       = 128
    "O2"       = 128
    "cchin_l"  = 256
    "ECG1_2"   = 128
    "spo2"     = 16
    "flow"     = 32
    "nas_pres" = 64
    "thorax"   = 32
    "abdomen"  = 32
    "snore"    = 256
    "lleg"     = 128
    "rleg"     = 128
    "position" = 16
    "cs_EEG"   = 100
    "cs_LOC"   = 100
    "cs_ROC"   = 100
    "cs_EMG"   = 100
    "cs_ECG"   = 100
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count

# Row 1 (B1:X1) holds the EDF file names: strip the ".edf" extension.
for ($c = 2; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value2
    if ($val -and $val.ToString().EndsWith(".edf")) {
        $cell.Value2 = $val.ToString().Substring(0, $val.ToString().Length - 4)
    }
}

# Rows 2..25 hold the channel name repeated across every file column;
# append ": <sample rate>" taken from the pyedflib read-out.
for ($r = 2; $r -le $lastRow; $r++) {
    $name = $ws.Cells.Item($r, 2).Value2
    if ($null -eq $name) { continue }
    $key = $name.ToString()
    if ($rates.ContainsKey($key)) {
        $newVal = "$key`: $($rates[$key])"
        for ($c = 2; $c -le $lastCol; $c++) {
            $ws.Cells.Item($r, $c).Value2 = $newVal
        }
    }
}

# Auto-fit the columns now that the text is wider: the index column (A)
# stays narrow, the 23 file-name columns (B:X) widen to fit "xxxx: 999".
$ws.Columns.AutoFit() | Out-Null
$ws.Range("A:A").ColumnWidth = 3
$ws.Range("B:X").ColumnWidth = 13.428571428571429

# Drop the explicit page setup / printer settings.
$ws.PageSetup.PrintArea = ""

# Move the active selection down, mirroring the user's post-edit click.
$ws.Range("B28").Select() | Out-Null
